# Applies the SOUTH_CAROLINA_2021.xlsx edit:
#  1. Remove the trailing footnote/metadata rows (old rows 1315-1319).
#  2. Rename the header row to short machine-friendly column names.
#  3. Title-case Spanish connector words ("de", "del", "la", "las", "los",
#     "el", "y") inside state/municipality names in columns A and B
#     (but never the first word of a cell).
#  4. Fix the one-off "MonteMorelos" typo to "Montemorelos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the trailing metadata rows -------------------------------
$ws.Range("A1315:A1319").EntireRow.Delete() | Out-Null

# --- 2. Rename header row ------------------------------------------------
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# --- 3. Title-case connector words in columns A and B (rows 2-1313) -----
$connectors = @("de", "del", "la", "las", "los", "el", "y")
$lastRow = 1313

for ($col = 1; $col -le 2; $col++) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value2
        if ($val -eq $null) { continue }
        if ($val.GetType().Name -ne "String") { continue }

        $words = $val -split " "
        $changed = $false
        for ($i = 1; $i -lt $words.Length; $i++) {
            $w = $words[$i]
            if ($connectors -contains $w.ToLower()) {
                $words[$i] = $w.Substring(0, 1).ToUpper() + $w.Substring(1).ToLower()
                $changed = $true
            }
        }

        if ($changed) {
            $cell.Value = ($words -join " ")
        }
    }
}

# --- 4. One-off capitalization fix ---------------------------------------
$b606 = $ws.Cells.Item(606, 2)
if ($b606.Value2 -eq "MonteMorelos") {
    $b606.Value = "Montemorelos"
}
